# Update cryptos list with latest prices and 1h volume change percentages
# (mirrors the scheduled GitHub Actions data refresh for cryptos.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.638.46"
$ws.Range("E2").Value = "  -1.63%  "

# Row 3
$ws.Range("D3").Value = "1.589.41"
$ws.Range("E3").Value = "  -2.07%  "

# Row 4
$ws.Range("E4").Value = "  +0.14%  "

# Row 5
$ws.Range("D5").Value = "'210.80"
$ws.Range("E5").Value = "  -1.59%  "

# Row 6
$ws.Range("D6").Value = "'0.509"
$ws.Range("E6").Value = "  -1.87%  "

# Row 7
$ws.Range("E7").Value = "  +0.15%  "

# Row 8
$ws.Range("E8").Value = "  -1.97%  "

# Row 9
$ws.Range("D9").Value = "'0.0615"
$ws.Range("E9").Value = "  -1.74%  "

# Row 10
$ws.Range("D10").Value = "'19.64"
$ws.Range("E10").Value = "  -3.17%  "

# Row 11
$ws.Range("D11").Value = "'0.0834"
$ws.Range("E11").Value = "  -1.29%  "

# Row 12
$ws.Range("D12").Value = "1.812.34"
$ws.Range("E12").Value = "  -2.05%  "

# Row 13
$ws.Range("D13").Value = "1.598.15"
$ws.Range("E13").Value = "  -1.49%  "

# Row 14
$ws.Range("D14").Value = "'4.02"
$ws.Range("E14").Value = "  -2.38%  "

# Row 15
$ws.Range("D15").Value = "'0.522"
$ws.Range("E15").Value = "  -3.65%  "

# Row 16
$ws.Range("D16").Value = "'64.83"
$ws.Range("E16").Value = "  +0.61%  "

# Row 17
$ws.Range("D17").Value = "26.624.60"
$ws.Range("E17").Value = "  -1.59%  "

# Row 18
$ws.Range("E18").Value = "  -2.17%  "

# Row 19
$ws.Range("D19").Value = "'208.51"
$ws.Range("E19").Value = "  -3.38%  "

# Row 20
$ws.Range("E20").Value = "  +0.15%  "

# Row 21
$ws.Range("D21").Value = "'6.73"
$ws.Range("E21").Value = "  -2.45%  "

# Row 22
$ws.Range("E22").Value = "  -3.05%  "

# Row 23
$ws.Range("D23").Value = "'2.33"
$ws.Range("E23").Value = "  -3.45%  "

# Row 24
$ws.Range("D24").Value = "'8.87"
$ws.Range("E24").Value = "  -1.69%  "

# Row 25
$ws.Range("D25").Value = "'146.94"
$ws.Range("E25").Value = "  -0.09%  "

# Row 26
$ws.Range("E26").Value = "  +0.21%  "

# Row 27
$ws.Range("E27").Value = "  -0.73%  "

# Row 28
$ws.Range("D28").Value = "'0.113"
$ws.Range("E28").Value = "  -3.33%  "

# Row 29
$ws.Range("D29").Value = "'15.30"
$ws.Range("E29").Value = "  -1.74%  "

# Row 30
$ws.Range("D30").Value = "'0.0507"
$ws.Range("E30").Value = "  +0.32%  "

# Row 31
$ws.Range("E31").Value = "  -1.66%  "

# Row 32
$ws.Range("D32").Value = "'3.22"
$ws.Range("E32").Value = "  -3.84%  "

# Row 33
$ws.Range("E33").Value = "  +21.36%  "

# Row 34
$ws.Range("E34").Value = "  -2.86%  "

# Row 35
$ws.Range("D35").Value = "1.300.77"
$ws.Range("E35").Value = "  -2.97%  "

# Row 36
$ws.Range("E36").Value = "  -1.18%  "

# Row 37
$ws.Range("D37").Value = "'1.48"
$ws.Range("E37").Value = "  -5.27%  "

# Row 38
$ws.Range("D38").Value = "'0.0171"
$ws.Range("E38").Value = "  -2.76%  "

# Row 39
$ws.Range("E39").Value = "  -2.31%  "

# Row 40
$ws.Range("E40").Value = "  +0.17%  "

# Row 41
$ws.Range("E41").Value = "  -1.29%  "

# Row 42
$ws.Range("D42").Value = "'5.36"
$ws.Range("E42").Value = "  +2.73%  "

# Row 43
$ws.Range("E43").Value = "  -2.81%  "

# Row 44
$ws.Range("D44").Value = "'62.72"
$ws.Range("E44").Value = "  -4.09%  "

# Row 45
$ws.Range("D45").Value = "1.725.43"
$ws.Range("E45").Value = "  -1.86%  "

# Row 46
$ws.Range("D46").Value = "'89.70"
$ws.Range("E46").Value = "  -0.78%  "

# Row 47
$ws.Range("D47").Value = "'1.60"
$ws.Range("E47").Value = "  -0.23%  "

# Row 48
$ws.Range("D48").Value = "'0.834"
$ws.Range("E48").Value = "  -2.41%  "

# Row 49
$ws.Range("E49").Value = "  -1.29%  "

# Row 50
$ws.Range("D50").Value = "'0.0503"
$ws.Range("E50").Value = "  -1.72%  "

# Row 51
$ws.Range("D51").Value = "'7.50"
$ws.Range("E51").Value = "  -0.86%  "
